# Apply the commit: "Updated trans/BVTStL so that all vehicle types except
# aircraft are subject to the LCFS"
#
# 1) On the "About" sheet: collapse the old 5-line note (A15:A19) about
#    exempting aircraft/rail/ships into a single shorter note on A15 that
#    only mentions aircraft, and remove the now-unused trailing rows.
# 2) On the "BVTStL" sheet: flip rail (row 5) and ships (row 6) from 0 to 1
#    for both the passenger (B) and freight (C) columns, so only aircraft
#    (row 4) remains 0/0.

$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")
$bvtstl = $wb.Worksheets.Item("BVTStL")

# --- About sheet: replace the long note with the shortened one, and drop
# the now-empty trailing rows (16-19).
$about.Range("A15").Value = "Based on the California LCFS, we choose to exempt aircraft."
$about.Range("A16").Value = $null
$about.Range("A17").Value = $null
$about.Range("A18").Value = $null
$about.Range("A19").Value = $null

# --- BVTStL sheet: rail and ships are now subject to the LCFS too.
$bvtstl.Range("B5").Value = 1
$bvtstl.Range("C5").Value = 1
$bvtstl.Range("B6").Value = 1
$bvtstl.Range("C6").Value = 1

# Leave the selections where the author left them when the workbook was saved.
$bvtstl.Activate()
$bvtstl.Range("C4").Select() | Out-Null
$about.Activate()
$about.Range("C18").Select() | Out-Null

$wb.Save()
